# This script applies corrections to the "Stage" worksheet:
#  - Column A (Numéro): date prefix corrected from 2017 to 2015
#  - Column E (Moyenne de l'étudiant): a number of grade values corrected
#
# Commit message: correction problèmes génération fichiers json

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stage")

$ws.Cells.Item(3, 1).Value = 20150926
$ws.Cells.Item(4, 1).Value = 20150927
$ws.Cells.Item(4, 5).Value = 7
$ws.Cells.Item(5, 1).Value = 20150928
$ws.Cells.Item(6, 1).Value = 20150929
$ws.Cells.Item(6, 5).Value = 20
$ws.Cells.Item(7, 1).Value = 20150930
$ws.Cells.Item(7, 5).Value = 10
$ws.Cells.Item(8, 1).Value = 20150931
$ws.Cells.Item(8, 5).Value = 19
$ws.Cells.Item(9, 1).Value = 20150932
$ws.Cells.Item(9, 5).Value = 17
$ws.Cells.Item(10, 1).Value = 20150933
$ws.Cells.Item(10, 5).Value = 18
$ws.Cells.Item(11, 1).Value = 20150934
$ws.Cells.Item(11, 5).Value = 14
$ws.Cells.Item(12, 1).Value = 20150935
$ws.Cells.Item(12, 5).Value = 16
$ws.Cells.Item(13, 1).Value = 20150936
$ws.Cells.Item(13, 5).Value = 12
$ws.Cells.Item(14, 1).Value = 20150937
$ws.Cells.Item(14, 5).Value = 15
$ws.Cells.Item(15, 1).Value = 20150938
$ws.Cells.Item(15, 5).Value = 17
$ws.Cells.Item(16, 1).Value = 20150939
$ws.Cells.Item(16, 5).Value = 18
$ws.Cells.Item(17, 1).Value = 20150940
$ws.Cells.Item(17, 5).Value = 19
$ws.Cells.Item(18, 1).Value = 20150941
$ws.Cells.Item(18, 5).Value = 16
$ws.Cells.Item(19, 1).Value = 20150942
$ws.Cells.Item(19, 5).Value = 19
$ws.Cells.Item(20, 1).Value = 20150943
$ws.Cells.Item(20, 5).Value = 19
$ws.Cells.Item(21, 1).Value = 20150944
$ws.Cells.Item(21, 5).Value = 18
$ws.Cells.Item(22, 1).Value = 20150945
$ws.Cells.Item(22, 5).Value = 8
$ws.Cells.Item(23, 1).Value = 20150946
$ws.Cells.Item(23, 5).Value = 6
$ws.Cells.Item(24, 1).Value = 20150947
$ws.Cells.Item(25, 1).Value = 20150948
$ws.Cells.Item(25, 5).Value = 8
$ws.Cells.Item(26, 1).Value = 20150949
$ws.Cells.Item(26, 5).Value = 10
$ws.Cells.Item(27, 1).Value = 20150950
$ws.Cells.Item(27, 5).Value = 15
$ws.Cells.Item(28, 1).Value = 20150951
$ws.Cells.Item(28, 5).Value = 17
$ws.Cells.Item(29, 1).Value = 20150952
$ws.Cells.Item(29, 5).Value = 12
$ws.Cells.Item(30, 1).Value = 20150953
$ws.Cells.Item(30, 5).Value = 13
$ws.Cells.Item(31, 1).Value = 20150954
$ws.Cells.Item(31, 5).Value = 9
$ws.Cells.Item(32, 1).Value = 20150955
$ws.Cells.Item(32, 5).Value = 20
$ws.Cells.Item(33, 1).Value = 20150956
$ws.Cells.Item(33, 5).Value = 12
$ws.Cells.Item(34, 1).Value = 20150957
$ws.Cells.Item(34, 5).Value = 19
$ws.Cells.Item(35, 1).Value = 20150958
$ws.Cells.Item(35, 5).Value = 16
$ws.Cells.Item(36, 1).Value = 20150959
$ws.Cells.Item(36, 5).Value = 17
$ws.Cells.Item(37, 1).Value = 20150960
$ws.Cells.Item(37, 5).Value = 9
$ws.Cells.Item(38, 1).Value = 20150961
$ws.Cells.Item(38, 5).Value = 11
$ws.Cells.Item(39, 1).Value = 20150962
$ws.Cells.Item(39, 5).Value = 18
$ws.Cells.Item(40, 1).Value = 20150963
$ws.Cells.Item(40, 5).Value = 12
$ws.Cells.Item(41, 1).Value = 20150964
$ws.Cells.Item(41, 5).Value = 10
$ws.Cells.Item(42, 1).Value = 20150965
$ws.Cells.Item(42, 5).Value = 13
$ws.Cells.Item(43, 1).Value = 20150966
$ws.Cells.Item(43, 5).Value = 7
$ws.Cells.Item(44, 1).Value = 20150967
$ws.Cells.Item(44, 5).Value = 10
$ws.Cells.Item(45, 1).Value = 20150968
$ws.Cells.Item(45, 5).Value = 5
$ws.Cells.Item(46, 1).Value = 20150969
$ws.Cells.Item(46, 5).Value = 9
$ws.Cells.Item(47, 1).Value = 20150970
$ws.Cells.Item(47, 5).Value = 19
$ws.Cells.Item(48, 1).Value = 20150971
$ws.Cells.Item(48, 5).Value = 19
$ws.Cells.Item(49, 1).Value = 20150972
$ws.Cells.Item(49, 5).Value = 11
$ws.Cells.Item(50, 1).Value = 20150973
$ws.Cells.Item(50, 5).Value = 11
$ws.Cells.Item(51, 1).Value = 20150974
$ws.Cells.Item(51, 5).Value = 5
$ws.Cells.Item(52, 1).Value = 20150975
$ws.Cells.Item(52, 5).Value = 12
$ws.Cells.Item(53, 1).Value = 20150976
$ws.Cells.Item(53, 5).Value = 8
$ws.Cells.Item(54, 1).Value = 20150977
$ws.Cells.Item(54, 5).Value = 12
$ws.Cells.Item(55, 1).Value = 20150978
$ws.Cells.Item(55, 5).Value = 17
$ws.Cells.Item(56, 1).Value = 20150979
$ws.Cells.Item(56, 5).Value = 5
$ws.Cells.Item(57, 1).Value = 20150980
$ws.Cells.Item(57, 5).Value = 16
$ws.Cells.Item(58, 1).Value = 20150981
$ws.Cells.Item(58, 5).Value = 5
$ws.Cells.Item(59, 1).Value = 20150982
$ws.Cells.Item(59, 5).Value = 7
$ws.Cells.Item(60, 1).Value = 20150983
$ws.Cells.Item(60, 5).Value = 6
$ws.Cells.Item(61, 1).Value = 20150984
$ws.Cells.Item(61, 5).Value = 6
$ws.Cells.Item(62, 1).Value = 20150985
$ws.Cells.Item(62, 5).Value = 17
$ws.Cells.Item(63, 1).Value = 20150986
$ws.Cells.Item(63, 5).Value = 11
